$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$np = $s.NotesPage
for ($i = $np.Shapes.Count; $i -ge 1; $i--) {
    $np.Shapes.Item($i).Delete()
}
Write-Output ("Shapes count after delete: " + $np.Shapes.Count)
